$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.271733999252319
$ws.Range("B1").Value = 1.86004102230072
$ws.Range("C1").Value = 5.628185749053955
$ws.Range("D1").Value = 1.939043521881104
$ws.Range("E1").Value = 1.116713047027588
